$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "weight" column (J) after the existing "order" column (I).
# Copy the header cell's formatting (bold/fill style) from I1 into J1, then
# set its text.
[void]$ws.Range("I1").Copy()
[void]$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "weight"

# Fill the weight values for all 15 data rows with 1 (default/no style,
# matching the other newly-added plain-numeric cells).
$ws.Range("J2:J16").Value = 1

# Update the active selection to match the edited workbook state.
[void]$ws.Range("K5").Select()
